$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.835.32"
$ws.Range("E2").Value = "  +5.57%  "
$ws.Range("D3").Value = "2.276.55"
$ws.Range("E3").Value = "  +3.75%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'232.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").Value = "'0.637"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.65%  "
$ws.Range("D7").Value = "'64.79"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.23%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.429"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.53%  "
$ws.Range("D10").Value = "'0.104"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +17.67%  "
$ws.Range("D11").Value = "'57.53"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("E12").Value = "  +18.83%  "
$ws.Range("D13").Value = "'0.103"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").Value = "2.611.52"
$ws.Range("E14").Value = "  +3.57%  "
$ws.Range("D15").Value = "'15.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.47%  "
$ws.Range("D16").Value = "'5.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.15%  "
$ws.Range("D17").Value = "'0.825"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.80%  "
$ws.Range("D18").Value = "2.269.67"
$ws.Range("E18").Value = "  +3.22%  "
$ws.Range("D19").Value = "43.714.17"
$ws.Range("E19").Value = "  +5.43%  "
$ws.Range("D20").Value = "'0.0000101"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +12.72%  "
$ws.Range("D21").Value = "'74.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.35%  "
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").Value = "'250.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.93%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").Value = "'2.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.01%  "
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("D27").Value = "'10.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.84%  "
$ws.Range("D28").Value = "'173.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.82%  "
$ws.Range("D29").Value = "'20.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.88%  "
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("E31").Value = "  +0.59%  "
$ws.Range("D32").Value = "'2.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.83%  "
$ws.Range("D33").Value = "'0.123"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.68%  "
$ws.Range("D34").Value = "'0.0690"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.37%  "
$ws.Range("D35").Value = "'5.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.22%  "
$ws.Range("D36").Value = "'4.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.54%  "
$ws.Range("E37").Value = "  +9.29%  "
$ws.Range("D38").Value = "'6.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.20%  "
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("E40").Value = "  +6.28%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "'17.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.68%  "
$ws.Range("D43").Value = "'8.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "'4.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.09%  "
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").Value = "'10.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +24.25%  "
$ws.Range("D46").Value = "'0.0971"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("E47").Value = "  +1.83%  "
$ws.Range("D48").Value = "'98.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("D49").Value = "1.480.68"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("E50").Value = "  +6.79%  "
$ws.Range("D51").Value = "'0.000204"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -14.29%  "
